$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.492.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.821.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.29"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3889"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08506"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +9.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.81"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.112"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.440"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.04"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.515"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.817.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001142"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.77"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06596"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.77"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.091"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.530.67"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.275"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.39"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.025.44"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.399"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.48"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.725"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07432"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.645"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02359"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.213"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.814"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6313"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.28"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.194"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.786"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5953"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "126.40"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06979"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.36"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.17%  "
